# Scheduled market-price refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets with the
# latest pulled values.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 498.66666
$ws.Range("I4").Value = 519
$ws.Range("J4").Value = 397
$ws.Range("K4").Value = 519
$ws.Range("L4").Value = 397
$ws.Range("M4").Value = -405
$ws.Range("N4").Value = -625
$ws.Range("H5").Value = 228.18182
$ws.Range("I5").Value = 83.28570999999999
$ws.Range("K5").Value = 83.28570999999999
$ws.Range("M5").Value = 31.71429000000001
$ws.Range("H12").Value = 2038.2307
$ws.Range("I12").Value = 1900
$ws.Range("K12").Value = 1900
$ws.Range("M12").Value = -1730
$ws.Range("H74").Value = 2375.75
$ws.Range("I74").Value = 1834.3334
$ws.Range("K74").Value = 1834.3334
$ws.Range("M74").Value = -898.3334
$ws.Range("H77").Value = 2375.75
$ws.Range("I77").Value = 1834.3334
$ws.Range("K77").Value = 9171.666999999999
$ws.Range("M77").Value = -4491.666999999999
$ws.Range("H88").Value = 2260.6924
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 2382.4167
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 2382.4167
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -3194.4167
$ws.Range("H91").Value = 2260.6924
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 2382.4167
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 2382.4167
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -5190.4167
$ws.Range("H111").Value = 2724.8
$ws.Range("I111").Value = 2341.6667
$ws.Range("K111").Value = 7025.000100000001
$ws.Range("M111").Value = -3958.000100000001
$ws.Range("H132").Value = 3799.4
$ws.Range("I132").Value = 2998.5
$ws.Range("K132").Value = 8995.5
$ws.Range("M132").Value = -6465.5
$ws.Range("H137").Value = 1553.5385
$ws.Range("I137").Value = 1381.4546
$ws.Range("K137").Value = 4144.3638
$ws.Range("M137").Value = -1594.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 54
$ws.Range("I4").Value = 48.5
$ws.Range("K4").Value = 48.5
$ws.Range("M4").Value = 67.5
$ws.Range("H5").Value = 716.5
$ws.Range("J5").Value = 1050
$ws.Range("L5").Value = 1050
$ws.Range("N5").Value = -1274
$ws.Range("H97").Value = 78.166664
$ws.Range("I97").Value = 98.666664
$ws.Range("J97").Value = 57.666668
$ws.Range("K97").Value = 98.666664
$ws.Range("L97").Value = 57.666668
$ws.Range("M97").Value = 397.333336
$ws.Range("N97").Value = -1049.666668
$ws.Range("H102").Value = 2308.9092
$ws.Range("I102").Value = 733.1667
$ws.Range("J102").Value = 4199.8
$ws.Range("K102").Value = 733.1667
$ws.Range("L102").Value = 4199.8
$ws.Range("M102").Value = 888.8333
$ws.Range("N102").Value = -7443.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 716.5
$ws.Range("J4").Value = 1050
$ws.Range("L4").Value = 1050
$ws.Range("N4").Value = -1280
$ws.Range("H94").Value = 2551.889
$ws.Range("J94").Value = 4200
$ws.Range("L94").Value = 4200
$ws.Range("N94").Value = -5102
$ws.Range("H99").Value = 2036
$ws.Range("I99").Value = 1883.5714
$ws.Range("K99").Value = 1883.5714
$ws.Range("M99").Value = -385.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2333.7368
$ws.Range("I7").Value = 1454.8
$ws.Range("K7").Value = 1454.8
$ws.Range("M7").Value = -1341.8
$ws.Range("H42").Value = 2000
$ws.Range("I42").Value = 2000
$ws.Range("K42").Value = 2000
$ws.Range("M42").Value = -1407

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 600
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -938
$ws.Range("H29").Value = 220.57143
$ws.Range("I29").Value = 185
$ws.Range("J29").Value = 226.5
$ws.Range("K29").Value = 555
$ws.Range("L29").Value = 679.5
$ws.Range("M29").Value = -278
$ws.Range("N29").Value = -1233.5
$ws.Range("H30").Value = 200
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 600
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -804
$ws.Range("H55").Value = 1010.1667
$ws.Range("I55").Value = 638.5
$ws.Range("J55").Value = 1084.5
$ws.Range("K55").Value = 1915.5
$ws.Range("L55").Value = 3253.5
$ws.Range("M55").Value = -1738.5
$ws.Range("N55").Value = -3607.5
$ws.Range("H58").Value = 4000
$ws.Range("H99").Value = 12000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1329
$ws.Range("I9").Value = 998.5
$ws.Range("J9").Value = 1990
$ws.Range("K9").Value = 998.5
$ws.Range("L9").Value = 1990
$ws.Range("M9").Value = -828.5
$ws.Range("N9").Value = -2330
$ws.Range("H70").Value = 24900
$ws.Range("I70").Value = 24900
$ws.Range("K70").Value = 24900
$ws.Range("M70").Value = -24630
$ws.Range("H73").Value = 24900
$ws.Range("I73").Value = 24900
$ws.Range("K73").Value = 24900
$ws.Range("M73").Value = -23964
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 421.52173
$ws.Range("I97").Value = 391.88235
$ws.Range("K97").Value = 391.88235
$ws.Range("M97").Value = 104.11765
$ws.Range("H105").Value = 42555
$ws.Range("J105").Value = 42555
$ws.Range("L105").Value = 42555
$ws.Range("N105").Value = -49543

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 40000
$ws.Range("J38").Value = 40000
$ws.Range("L38").Value = 40000
$ws.Range("N38").Value = -40820
$ws.Range("H40").Value = 1855
$ws.Range("I40").Value = 1855
$ws.Range("K40").Value = 1855
$ws.Range("M40").Value = -1719
$ws.Range("H46").Value = 2842.4211
$ws.Range("I46").Value = 2363.6365
$ws.Range("J46").Value = 3500.75
$ws.Range("K46").Value = 2363.6365
$ws.Range("L46").Value = 3500.75
$ws.Range("M46").Value = -2175.6365
$ws.Range("N46").Value = -3876.75
$ws.Range("H58").Value = 11401.667
$ws.Range("I58").Value = 8000
$ws.Range("J58").Value = 13102.5
$ws.Range("K58").Value = 8000
$ws.Range("L58").Value = 13102.5
$ws.Range("M58").Value = -7740
$ws.Range("N58").Value = -13622.5
$ws.Range("H122").Value = 5311.2144
$ws.Range("I122").Value = 5027.4614
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 15082.3842
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -12632.3842
$ws.Range("N122").Value = -31900
$ws.Range("H136").Value = 4154.5557
$ws.Range("I136").Value = 3731.8333
$ws.Range("K136").Value = 11195.4999
$ws.Range("M136").Value = -8645.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 2000
$ws.Range("K21").Value = 2000
$ws.Range("M21").Value = -1765
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1710
$ws.Range("H126").Value = 974.9091
$ws.Range("I126").Value = 1025
$ws.Range("J126").Value = 749.5
$ws.Range("K126").Value = 3075
$ws.Range("L126").Value = 2248.5
$ws.Range("M126").Value = -605
$ws.Range("N126").Value = -7188.5

Write-Output "Updated 209 cells across 8 profession sheets."
